# Exterior Lighting Control design doc update
# Applies the changes described by the commit "Exterior Lighting Control design doc update"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: "Enable Advanced RTU Controls" -> "Exterior Lighting Control"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Enable Advanced RTU Controls", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Exterior Lighting Control", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Description paragraph: merge "educe" + "s" + " all exterior..." runs
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("educe" + "s" + " all exterior lighting to 30% of its peak power between midnight or within 1 hour of business closing, whichever is later, and until 6 am or business opening, whichever is earlier, and during any period activity is not detected for a time longer than 15 minutes.", `
                         $false, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "educes all exterior lighting to 30% of its peak power between midnight or within 1 hour of business closing, whichever is later, and until 6 am or business opening, whichever is earlier, and during any period activity is not detected for a time longer than 15 minutes.", `
                         2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Modeler Description paragraph: rewrite
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("This EEM loops through all occupancy schedules in the model, determines the earliest and latest building open and close times, and creates a new fractional schedule for exterior lights based on these times. The new schedule reduces the exterior lighting power to from 1.0 to 0.7 (30% reduction) during this interval.", `
                         $false, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "This measure first loops through all occupancy schedules in the model, determines the earliest and latest occupied, compares these times to a fixed 0000-0600 schedule, and creates a new fractional schedule for exterior lights with the shortest interval. The new schedule reduces the exterior lighting power to from 1.0 to 0.7 (30% reduction) during this interval. The measure then loops through all exterior lights objects in the model, changes the control option to ScheduleNameOnly if necessary, and sets the schedule to the new exterior lights schedule.", `
                         2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Initial Condition Message paragraph
# ---------------------------------------------------------------------------
$pInitial = $d.Paragraphs(13)
$pInitial.Range.Find.Execute("Total exterior lighting power = ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Number of exterior lights objects in model = ", 2) | Out-Null
$pInitial.Range.Find.Execute("TODO} W  ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "ext_ltg_count}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Final Condition Message paragraph
# ---------------------------------------------------------------------------
$pFinal = $d.Paragraphs(15)
$pFinal.Range.Find.Execute("Total exterior lighting power = ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Number of exterior lights objects changed = ", 2) | Out-Null
$pFinal.Range.Find.Execute("TODO} W", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "ext_ltg_changed}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Insert new "Info Messages" section after Final Condition Message, before
#    "Not Applicable Messages"
# ---------------------------------------------------------------------------
$pFinal.Range.InsertParagraphAfter()
$pInfoHeading = $d.Paragraphs($pFinal.Index + 1)
$pInfoHeading.Style = "Heading3"
$pInfoHeading.Range.Text = "Info Messages"

$prev = $pInfoHeading
$infoLines = @(
  "Adding new schedule to model: #{ext_lights_sch_name}",
  "Applying exterior lighting controls to: #{el.name}",
  "=> control option set to: AstronomicalClock",
  "=> setting control option to: ScheduleNameOnly",
  "=> setting schedule to: #{ext_lights_sch.name}"
)
foreach ($line in $infoLines) {
    $prev.Range.InsertParagraphAfter()
    $np = $d.Paragraphs($prev.Index + 1)
    $np.Style = "Normal"
    $np.SpaceAfter = 0
    $np.Range.Text = $line
    $prev = $np
}

# ---------------------------------------------------------------------------
# 7. Not Applicable Messages paragraph: merge the 3 runs into 1
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("No exterior lighting found, EEM not applied.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "No exterior lights objects found in model.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8. Warning Messages content paragraph -> bookmark + "NA"
#    (the bookmark moves here from the end of the document)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Exterior lighting control option changed from Astronomical Clock to Scheduled.", `
                         $false, $false, $false, $false, $false, `
                         $true, 1, $false, "NA", 2) | Out-Null

$pWarning = $d.Paragraphs(19)
$bmRange = $pWarning.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 9. "Test results: TBD" paragraph -> replaced by new Tests block
# ---------------------------------------------------------------------------
$pTests = $d.Paragraphs($d.Paragraphs.Count - 1)

$pTests.Range.Text = ""
$pTests.Style = "Normal"
$pTests.SpaceAfter = 0

$prev = $pTests
$prev.Range.InsertParagraphAfter()
$pVerify = $d.Paragraphs($prev.Index + 1)
$pVerify.Style = "Normal"
$pVerify.SpaceAfter = 0
$pVerify.Range.Text = "The following tests should be verified:"
$prev = $pVerify

$prev.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs($prev.Index + 1)
$pBlank.Style = "Normal"
$pBlank.SpaceAfter = 0
$prev = $pBlank

$bulletLines0 = @(
  "A new Schedule Ruleset is created for exterior lights ",
  "Exterior Lights Schedule is set to 0.7 from 0000 ` 0600",
  "All Exterior Lights objects in the model are changed to:"
)

# first list item needs its own list template reference: reuse the existing
# numId 9 list by cloning an existing ilvl-0 item from the Use Case Types list
$pTemplate0 = $d.Paragraphs(9)
foreach ($line in $bulletLines0) {
    $prev.Range.InsertParagraphAfter()
    $np = $d.Paragraphs($prev.Index + 1)
    $np.Style = "ListParagraph"
    $np.Range.ListFormat.ListTemplate = $pTemplate0.Range.ListFormat.ListTemplate
    $np.SpaceAfter = 0
    $np.Range.Text = $line
    $prev = $np
}

$bulletLines1 = @(
  "Schedule = Exterior Lights Schedule",
  "Control Option = ScheduleNameOnly"
)
foreach ($line in $bulletLines1) {
    $prev.Range.InsertParagraphAfter()
    $np = $d.Paragraphs($prev.Index + 1)
    $np.Style = "ListParagraph"
    $np.Range.ListFormat.ListTemplate = $pTemplate0.Range.ListFormat.ListTemplate
    $np.Range.ListFormat.ListIndent()
    $np.SpaceAfter = 0
    $np.Range.Text = $line
    $prev = $np
}

Write-Output "DONE"
